$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows 20-21 have the same date number format as column D elsewhere
$dateFmt = $ws.Range("D19").NumberFormat
$ws.Range("D20:D21").NumberFormat = $dateFmt

# Fill constant columns for newly added rows 20 and 21
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103003
$ws.Range("J20").Value = "Damasco"
$ws.Range("K20").Value = "Castle Brite"

$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103003
$ws.Range("J21").Value = "Damasco"
$ws.Range("K21").Value = "Castle Brite"

# Update the data columns (D, L-T) for every data row, 2 through 21
$ws.Range("D2").Value = 44547
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 350
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("Q2").Value = "`$/caja 18 kilos"
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1111
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 44547
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 350
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = "`$/caja 18 kilos"
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44547
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 350
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "`$/caja 18 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 889
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44159
$ws.Range("L5").Value = "Tercera"
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 15500
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15750
$ws.Range("Q5").Value = "`$/caja 15 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1050
$ws.Range("T5").Value = 15

$ws.Range("D6").Value = 44169
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 500
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "`$/caja 15 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1033
$ws.Range("T6").Value = 15

$ws.Range("D7").Value = 44544
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 600
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19000
$ws.Range("Q7").Value = "`$/caja 18 kilos"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1056
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44544
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range("Q8").Value = "`$/caja 18 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 889
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44530
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = "`$/caja 18 kilos"
$ws.Range("R9").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S9").Value = 1139
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44537
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 500
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21000
$ws.Range("Q10").Value = "`$/caja 18 kilos"
$ws.Range("R10").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S10").Value = 1167
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44537
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 17000
$ws.Range("Q11").Value = "`$/caja 18 kilos"
$ws.Range("R11").Value = "Región del Maule"
$ws.Range("S11").Value = 944
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44162
$ws.Range("L12").Value = "Tercera"
$ws.Range("M12").Value = 500
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15500
$ws.Range("Q12").Value = "`$/caja 15 kilos"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1033
$ws.Range("T12").Value = 15

$ws.Range("D13").Value = 44166
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 600
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 16500
$ws.Range("Q13").Value = "`$/caja 15 kilos"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 1100
$ws.Range("T13").Value = 15

$ws.Range("D14").Value = 44533
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 350
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 24000
$ws.Range("P14").Value = 24000
$ws.Range("Q14").Value = "`$/caja 18 kilos"
$ws.Range("R14").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S14").Value = 1333
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44533
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 350
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = "`$/caja 18 kilos"
$ws.Range("R15").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S15").Value = 1111
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44533
$ws.Range("L16").Value = "Tercera"
$ws.Range("M16").Value = 350
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 17000
$ws.Range("Q16").Value = "`$/caja 18 kilos"
$ws.Range("R16").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S16").Value = 944
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 44176
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("Q17").Value = "`$/caja 15 kilos"
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 1033
$ws.Range("T17").Value = 15

$ws.Range("D18").Value = 44187
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 350
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("Q18").Value = "`$/caja 15 kilos"
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1067
$ws.Range("T18").Value = 15

$ws.Range("D19").Value = 44187
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 13000
$ws.Range("O19").Value = 13000
$ws.Range("P19").Value = 13000
$ws.Range("Q19").Value = "`$/caja 15 kilos"
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("S19").Value = 867
$ws.Range("T19").Value = 15

$ws.Range("D20").Value = 44194
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 15500
$ws.Range("Q20").Value = "`$/caja 15 kilos"
$ws.Range("R20").Value = "Región Metropolitana"
$ws.Range("S20").Value = 1033
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 44540
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 600
$ws.Range("N21").Value = 16000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 16000
$ws.Range("Q21").Value = "`$/caja 18 kilos"
$ws.Range("R21").Value = "Región del Maule"
$ws.Range("S21").Value = 889
$ws.Range("T21").Value = 18

Write-Host "Done"